$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "95.794.91"
$ws.Range("E2").Value = "  -0.25%  "
$ws.Range("D3").Value = "3.588.23"
$ws.Range("E3").Value = "  -2.39%  "
$ws.Range("E4").Value = "  +0.06%  "
$ws.Range("D5").Value = "238.02"
$ws.Range("E5").Value = "  -2.02%  "
$ws.Range("D6").Value = "654.48"
$ws.Range("E6").Value = "  +1.33%  "
$ws.Range("E7").Value = "  +3.45%  "
$ws.Range("E8").Value = "  -0.05%  "
$ws.Range("E9").Value = "  +0.11%  "
$ws.Range("E10").Value = "  +1.52%  "
$ws.Range("D11").Value = "3.593.23"
$ws.Range("E11").Value = "  -2.24%  "
$ws.Range("B12").Value = "TRON"
$ws.Range("C12").Value = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
$ws.Range("D12").Value = "0.203"
$ws.Range("E12").Value = "  +1.06%  "
$ws.Range("B13").Value = "Avalanche"
$ws.Range("C13").Value = "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"
$ws.Range("D13").Value = "42.96"
$ws.Range("E13").Value = "  -2.50%  "
$ws.Range("D14").Value = "6.46"
$ws.Range("E14").Value = "  +0.90%  "
$ws.Range("D15").Value = "4.257.17"
$ws.Range("E15").Value = "  -2.67%  "
$ws.Range("D16").Value = "95.684.32"
$ws.Range("E16").Value = "  -0.19%  "
$ws.Range("D17").Value = "0.0000255"
$ws.Range("E17").Value = "  -1.00%  "
$ws.Range("D18").Value = "3.582.26"
$ws.Range("E18").Value = "  -2.60%  "
$ws.Range("D19").Value = "12.77"
$ws.Range("E19").Value = "  -5.06%  "
$ws.Range("D20").Value = "7.76"
$ws.Range("E20").Value = "  -3.80%  "
$ws.Range("D21").Value = "18.01"
$ws.Range("E21").Value = "  -3.91%  "
$ws.Range("D22").Value = "0.495"
$ws.Range("E22").Value = "  +2.52%  "
$ws.Range("D23").Value = "3.46"
$ws.Range("E23").Value = "  +0.33%  "
$ws.Range("D24").Value = "511.88"
$ws.Range("E24").Value = "  -1.88%  "
$ws.Range("D25").Value = "7.06"
$ws.Range("E25").Value = "  +3.35%  "
$ws.Range("E26").Value = "  -0.06%  "
$ws.Range("D27").Value = "96.06"
$ws.Range("E27").Value = "  -1.51%  "
$ws.Range("D28").Value = "12.83"
$ws.Range("E28").Value = "  +1.01%  "
$ws.Range("D29").Value = "3.781.39"
$ws.Range("E29").Value = "  -2.34%  "
$ws.Range("E30").Value = "  -3.59%  "
$ws.Range("E31").Value = "  +2.56%  "
$ws.Range("D32").Value = "11.59"
$ws.Range("E32").Value = "  -1.33%  "
$ws.Range("E33").Value = "  +0.16%  "
$ws.Range("D34").Value = "0.997"
$ws.Range("E34").Value = "  +0.16%  "
$ws.Range("E35").Value = "  -1.22%  "
$ws.Range("D36").Value = "32.01"
$ws.Range("E36").Value = "  -3.54%  "
$ws.Range("E37").Value = "  +12.28%  "
$ws.Range("D38").Value = "0.565"
$ws.Range("E38").Value = "  -2.75%  "
$ws.Range("D39").Value = "8.61"
$ws.Range("E39").Value = "  +7.88%  "
$ws.Range("D40").Value = "598.73"
$ws.Range("E40").Value = "  +6.28%  "
$ws.Range("E41").Value = "  -0.99%  "
$ws.Range("E42").Value = "  +0.10%  "
$ws.Range("D43").Value = "1.88"
$ws.Range("E43").Value = "  +6.78%  "
$ws.Range("D44").Value = "0.916"
$ws.Range("E44").Value = "  -5.55%  "
$ws.Range("E45").Value = "  -0.12%  "
$ws.Range("E46").Value = "  +3.67%  "
$ws.Range("D47").Value = "34.47"
$ws.Range("E47").Value = "  +1.11%  "
$ws.Range("B48").Value = "WhiteBITCoin"
$ws.Range("C48").Value = "https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt"
$ws.Range("D48").Value = "23.46"
$ws.Range("E48").Value = "  -1.04%  "
$ws.Range("B49").Value = "VeChain"
$ws.Range("C49").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D49").Value = "0.0418"
$ws.Range("E49").Value = "  -2.99%  "
$ws.Range("E50").Value = "  -0.18%  "
$ws.Range("D51").Value = "8.26"
$ws.Range("E51").Value = "  -0.89%  "
